$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H53").Value = 269.6875
$ws.Range("I53").Value = 265.6
$ws.Range("J53").Value = 276.5
$ws.Range("K53").Value = 265.6
$ws.Range("L53").Value = 276.5
$ws.Range("M53").Value = 371.4
$ws.Range("N53").Value = -1550.5

$ws.Range("H62").Value = 12126
$ws.Range("I62").Value = 10502.5
$ws.Range("K62").Value = 10502.5
$ws.Range("M62").Value = -9878.5

$ws.Range("H65").Value = 12126
$ws.Range("I65").Value = 10502.5
$ws.Range("K65").Value = 52512.5
$ws.Range("M65").Value = -49392.5

$ws.Range("H74").Value = 16938.615
$ws.Range("I74").Value = 19874.143
$ws.Range("J74").Value = 4609.4
$ws.Range("K74").Value = 19874.143
$ws.Range("L74").Value = 4609.4
$ws.Range("M74").Value = -18938.143
$ws.Range("N74").Value = -6481.4

$ws.Range("H77").Value = 16938.615
$ws.Range("I77").Value = 19874.143
$ws.Range("J77").Value = 4609.4
$ws.Range("K77").Value = 99370.715
$ws.Range("L77").Value = 23047
$ws.Range("M77").Value = -94690.715
$ws.Range("N77").Value = -32407

$ws.Range("H100").Value = 3176.1667
$ws.Range("I100").Value = 2398
$ws.Range("K100").Value = 2398
$ws.Range("M100").Value = -1857

$ws.Range("H112").Value = 6023
$ws.Range("I112").Value = 2424.6667
$ws.Range("J112").Value = 7372.375
$ws.Range("K112").Value = 7274.000100000001
$ws.Range("L112").Value = 22117.125
$ws.Range("M112").Value = -6166.000100000001
$ws.Range("N112").Value = -24333.125

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 3228.0833
$ws.Range("I45").Value = 2637.6667
$ws.Range("K45").Value = 2637.6667
$ws.Range("M45").Value = -2260.6667

$ws.Range("H61").Value = 4365.8335
$ws.Range("I61").Value = 3914.8965
$ws.Range("J61").Value = 6234
$ws.Range("K61").Value = 3914.8965
$ws.Range("L61").Value = 6234
$ws.Range("M61").Value = -3702.8965
$ws.Range("N61").Value = -6658

$ws.Range("H74").Value = 4599.773
$ws.Range("I74").Value = 1981.7273
$ws.Range("K74").Value = 1981.7273
$ws.Range("M74").Value = -1107.7273

$ws.Range("H77").Value = 4599.773
$ws.Range("I77").Value = 1981.7273
$ws.Range("K77").Value = 9908.636500000001
$ws.Range("M77").Value = -5540.636500000001

$ws.Range("H136").Value = 4365.8335
$ws.Range("I136").Value = 3914.8965
$ws.Range("J136").Value = 6234
$ws.Range("K136").Value = 11744.6895
$ws.Range("L136").Value = 18702
$ws.Range("M136").Value = -9194.6895
$ws.Range("N136").Value = -23802

$ws.Range("H138").Value = 70000
$ws.Range("J138").Value = 70000
$ws.Range("L138").Value = 70000
$ws.Range("N138").Value = -80280

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 2589.0688
$ws.Range("I86").Value = 1041.1154
$ws.Range("J86").Value = 16004.667
$ws.Range("K86").Value = 1041.1154
$ws.Range("L86").Value = 16004.667
$ws.Range("M86").Value = 81.88460000000009
$ws.Range("N86").Value = -18250.667

$ws.Range("H89").Value = 2589.0688
$ws.Range("I89").Value = 1041.1154
$ws.Range("J89").Value = 16004.667
$ws.Range("K89").Value = 5205.576999999999
$ws.Range("L89").Value = 80023.33499999999
$ws.Range("M89").Value = 410.4230000000007
$ws.Range("N89").Value = -91255.33499999999

$ws.Range("H99").Value = 537.25
$ws.Range("I99").Value = 525
$ws.Range("J99").Value = 549.5
$ws.Range("K99").Value = 525
$ws.Range("L99").Value = 549.5
$ws.Range("M99").Value = 973
$ws.Range("N99").Value = -3545.5

$ws.Range("H107").Value = 1616.091
$ws.Range("I107").Value = 1222
$ws.Range("J107").Value = 2667
$ws.Range("K107").Value = 1222
$ws.Range("L107").Value = 2667
$ws.Range("M107").Value = 698
$ws.Range("N107").Value = -6507

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 1077
$ws.Range("I16").Value = 919.6
$ws.Range("J16").Value = 1601.6666
$ws.Range("K16").Value = 919.6
$ws.Range("L16").Value = 1601.6666
$ws.Range("M16").Value = -632.6
$ws.Range("N16").Value = -2175.6666

$ws.Range("H41").Value = 30000.5
$ws.Range("J41").Value = 30000.5
$ws.Range("L41").Value = 30000.5
$ws.Range("N41").Value = -30856.5

$ws.Range("H58").Value = 3092.3
$ws.Range("I58").Value = 2785.6
$ws.Range("J58").Value = 3399
$ws.Range("K58").Value = 2785.6
$ws.Range("L58").Value = 3399
$ws.Range("M58").Value = -2582.6
$ws.Range("N58").Value = -3805

$ws.Range("H96").Value = 46084.6
$ws.Range("J96").Value = 46084.6
$ws.Range("L96").Value = 46084.6
$ws.Range("N96").Value = -51576.6

$ws.Range("H99").Value = 5786.7
$ws.Range("I99").Value = 6385.6665
$ws.Range("J99").Value = 4888.25
$ws.Range("K99").Value = 6385.6665
$ws.Range("L99").Value = 4888.25
$ws.Range("M99").Value = -4887.6665
$ws.Range("N99").Value = -7884.25

$ws.Range("H113").Value = 1077
$ws.Range("I113").Value = 919.6
$ws.Range("J113").Value = 1601.6666
$ws.Range("K113").Value = 919.6
$ws.Range("L113").Value = 1601.6666
$ws.Range("M113").Value = 1250.4
$ws.Range("N113").Value = -5941.6666

$ws.Range("H126").Value = 5786.7
$ws.Range("I126").Value = 6385.6665
$ws.Range("J126").Value = 4888.25
$ws.Range("K126").Value = 19156.9995
$ws.Range("L126").Value = 14664.75
$ws.Range("M126").Value = -16686.9995
$ws.Range("N126").Value = -19604.75

$ws.Range("H136").Value = 3092.3
$ws.Range("I136").Value = 2785.6
$ws.Range("J136").Value = 3399
$ws.Range("K136").Value = 8356.799999999999
$ws.Range("L136").Value = 10197
$ws.Range("M136").Value = -5806.799999999999
$ws.Range("N136").Value = -15297

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 385.7143
$ws.Range("I2").Value = 200
$ws.Range("K2").Value = 1200
$ws.Range("M2").Value = -1087

$ws.Range("H12").Value = 155
$ws.Range("J12").Value = 155
$ws.Range("L12").Value = 465
$ws.Range("N12").Value = -811

$ws.Range("H17").Value = 982.375
$ws.Range("I17").Value = 982.375
$ws.Range("J17").Value = 0
$ws.Range("K17").Value = 2947.125
$ws.Range("L17").Value = 0
$ws.Range("M17").ClearContents()

$ws.Range("H34").Value = 17839.2
$ws.Range("J34").Value = 17839.2
$ws.Range("L34").Value = 53517.60000000001
$ws.Range("N34").Value = -53685.60000000001

$ws.Range("H39").Value = 12499.5
$ws.Range("J39").Value = 19999
$ws.Range("L39").Value = 59997
$ws.Range("N39").Value = -60585

$ws.Range("H55").Value = 13998.5
$ws.Range("J55").Value = 13998.5
$ws.Range("L55").Value = 41995.5
$ws.Range("N55").Value = -42349.5

$ws.Range("H70").Value = 20000
$ws.Range("J70").Value = 20000
$ws.Range("L70").Value = 60000
$ws.Range("N70").Value = -60630

$ws.Range("H73").Value = 20000
$ws.Range("J73").Value = 20000
$ws.Range("L73").Value = 60000
$ws.Range("N73").Value = -62184

$ws.Range("H113").Value = 685.95654
$ws.Range("I113").Value = 378.9
$ws.Range("J113").Value = 922.1539
$ws.Range("K113").Value = 1136.7
$ws.Range("L113").Value = 2766.4617
$ws.Range("M113").Value = 1033.3
$ws.Range("N113").Value = -7106.4617

$ws.Range("H133").Value = 9112.809999999999
$ws.Range("I133").Value = 4108.625
$ws.Range("J133").Value = 12192.308
$ws.Range("K133").Value = 12325.875
$ws.Range("L133").Value = 36576.924
$ws.Range("M133").Value = -7265.875
$ws.Range("N133").Value = -46696.924

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 3067.0557
$ws.Range("I46").Value = 1020
$ws.Range("J46").Value = 3651.9285
$ws.Range("K46").Value = 1020
$ws.Range("L46").Value = 3651.9285
$ws.Range("M46").Value = -832
$ws.Range("N46").Value = -4027.9285

$ws.Range("H61").Value = 4545.6924
$ws.Range("I61").Value = 1909.6
$ws.Range("K61").Value = 1909.6
$ws.Range("M61").Value = -1707.6

$ws.Range("H93").Value = 1745
$ws.Range("J93").Value = 3954.6
$ws.Range("L93").Value = 3954.6
$ws.Range("N93").Value = -6450.6

$ws.Range("H113").Value = 4545.6924
$ws.Range("I113").Value = 1909.6
$ws.Range("K113").Value = 1909.6
$ws.Range("M113").Value = 260.4000000000001

$ws.Range("H136").Value = 6607.3086
$ws.Range("I136").Value = 5760.32
$ws.Range("K136").Value = 17280.96
$ws.Range("M136").Value = -14730.96
